$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 84; existing rows 84-92 shift down to 86-94.
$ws.Rows("84:85").Insert()

# New row 84 data
$ws.Cells.Item(84, 1).Value = 5
$ws.Cells.Item(84, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(84, 3).Value = "Maule"
$ws.Cells.Item(84, 4).Value = 44826
$ws.Cells.Item(84, 5).Value = 7
$ws.Cells.Item(84, 6).Value = "Fruta"
$ws.Cells.Item(84, 7).Value = 100107
$ws.Cells.Item(84, 8).Value = "Otros"
$ws.Cells.Item(84, 9).Value = 100107002
$ws.Cells.Item(84, 10).Value = "Chirimoya"
$ws.Cells.Item(84, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(84, 12).Value = "Especial"
$ws.Cells.Item(84, 13).Value = 30
$ws.Cells.Item(84, 14).Value = 30000
$ws.Cells.Item(84, 15).Value = 30000
$ws.Cells.Item(84, 16).Value = 30000
$ws.Cells.Item(84, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(84, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(84, 19).Value = 3000
$ws.Cells.Item(84, 20).Value = 10

# New row 85 data
$ws.Cells.Item(85, 1).Value = 5
$ws.Cells.Item(85, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(85, 3).Value = "Maule"
$ws.Cells.Item(85, 4).Value = 44826
$ws.Cells.Item(85, 5).Value = 7
$ws.Cells.Item(85, 6).Value = "Fruta"
$ws.Cells.Item(85, 7).Value = 100107
$ws.Cells.Item(85, 8).Value = "Otros"
$ws.Cells.Item(85, 9).Value = 100107002
$ws.Cells.Item(85, 10).Value = "Chirimoya"
$ws.Cells.Item(85, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(85, 12).Value = "Segunda"
$ws.Cells.Item(85, 13).Value = 20
$ws.Cells.Item(85, 14).Value = 25000
$ws.Cells.Item(85, 15).Value = 25000
$ws.Cells.Item(85, 16).Value = 25000
$ws.Cells.Item(85, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(85, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(85, 19).Value = 2500
$ws.Cells.Item(85, 20).Value = 10

# Ensure the date column keeps its datetime number format on the two new rows
$ws.Range("D84:D85").NumberFormat = "YYYY-MM-DD HH:MM:SS"
